$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 304.5
$ws.Range("I52").Value = 304.5
$ws.Range("K52").Value = 913.5
$ws.Range("M52").Value = -753.5
$ws.Range("H135").Value = 1444.0975
$ws.Range("I135").Value = 949.0909
$ws.Range("K135").Value = 8541.8181
$ws.Range("M135").Value = -6006.8181

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 32883.695
$ws.Range("I132").Value = 33661.312
$ws.Range("K132").Value = 100983.936
$ws.Range("M132").Value = -98453.93599999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1791
$ws.Range("I20").Value = 1534.8
$ws.Range("J20").Value = 2157
$ws.Range("K20").Value = 1534.8
$ws.Range("L20").Value = 2157
$ws.Range("M20").Value = -1287.8
$ws.Range("N20").Value = -2651

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 9993.6
$ws.Range("J9").Value = 9993.6
$ws.Range("L9").Value = 9993.6
$ws.Range("N9").Value = -10329.6
$ws.Range("H31").Value = 2877.4348
$ws.Range("I31").Value = 2877.4348
$ws.Range("K31").Value = 2877.4348
$ws.Range("M31").Value = -2582.4348
$ws.Range("H34").Value = 2877.4348
$ws.Range("I34").Value = 2877.4348
$ws.Range("K34").Value = 2877.4348
$ws.Range("M34").Value = -2675.4348
$ws.Range("H132").Value = 4139.353
$ws.Range("I132").Value = 4210.5625
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 12631.6875
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -10101.6875
$ws.Range("N132").Value = -14060
$ws.Range("H134").Value = 35491.742
$ws.Range("I134").Value = 40288.668
$ws.Range("K134").Value = 120866.004
$ws.Range("M134").Value = -118331.004

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 780.3333
$ws.Range("I9").Value = 780.3333
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 2340.9999
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -2116.9999
$ws.Range("N9").Value = $null
$ws.Range("H16").Value = 2199.25
$ws.Range("I16").Value = 2199.25
$ws.Range("K16").Value = 6597.75
$ws.Range("M16").Value = -6424.75
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").Value = $null
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").Value = $null
$ws.Range("H25").Value = 7071.143
$ws.Range("I25").Value = 6449.5
$ws.Range("J25").Value = 7900
$ws.Range("K25").Value = 19348.5
$ws.Range("L25").Value = 23700
$ws.Range("M25").Value = -19179.5
$ws.Range("N25").Value = -24038
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").Value = $null
$ws.Range("H30").Value = 7071.143
$ws.Range("I30").Value = 6449.5
$ws.Range("J30").Value = 7900
$ws.Range("K30").Value = 19348.5
$ws.Range("L30").Value = 23700
$ws.Range("M30").Value = -19246.5
$ws.Range("N30").Value = -23904
$ws.Range("H34").Value = 2553.4666
$ws.Range("I34").Value = 2043.1428
$ws.Range("J34").Value = 3000
$ws.Range("K34").Value = 6129.428400000001
$ws.Range("L34").Value = 9000
$ws.Range("M34").Value = -6045.428400000001
$ws.Range("N34").Value = -9168
$ws.Range("H35").Value = 3722.6667
$ws.Range("I35").Value = 2251
$ws.Range("K35").Value = 6753
$ws.Range("M35").Value = -6465
$ws.Range("H40").Value = 244.28572
$ws.Range("J40").Value = 250
$ws.Range("L40").Value = 1000
$ws.Range("N40").Value = -1138
$ws.Range("H46").Value = 999.5
$ws.Range("J46").Value = 1999
$ws.Range("L46").Value = 5997
$ws.Range("N46").Value = -6179
$ws.Range("H48").Value = 3538.3845
$ws.Range("I48").Value = 2000
$ws.Range("J48").Value = 3666.5833
$ws.Range("K48").Value = 6000
$ws.Range("L48").Value = 10999.7499
$ws.Range("M48").Value = -5750
$ws.Range("N48").Value = -11499.7499
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = $null
$ws.Range("N49").Value = $null
$ws.Range("H57").Value = 13249
$ws.Range("I57").Value = 1997
$ws.Range("J57").Value = 16999.666
$ws.Range("K57").Value = 5991
$ws.Range("L57").Value = 50998.99800000001
$ws.Range("M57").Value = -5432
$ws.Range("N57").Value = -52116.99800000001
$ws.Range("H59").Value = 1467
$ws.Range("I59").Value = 1467
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 4401
$ws.Range("L59").Value = 0
$ws.Range("M59").Value = -3861
$ws.Range("N59").Value = $null
$ws.Range("H140").Value = 2714.611
$ws.Range("I140").Value = 2714.611
$ws.Range("K140").Value = 8143.833
$ws.Range("M140").Value = -2963.833

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 14817.066
$ws.Range("J7").Value = 4993.6665
$ws.Range("L7").Value = 4993.6665
$ws.Range("N7").Value = -5217.6665
$ws.Range("H40").Value = 2998.6924
$ws.Range("I40").Value = 1477.5555
$ws.Range("K40").Value = 1477.5555
$ws.Range("M40").Value = -1341.5555
$ws.Range("H126").Value = 14817.066
$ws.Range("J126").Value = 4993.6665
$ws.Range("L126").Value = 14980.9995
$ws.Range("N126").Value = -19920.9995

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 9199.4
$ws.Range("I81").Value = 3349.5
$ws.Range("J81").Value = 13099.333
$ws.Range("K81").Value = 6699
$ws.Range("L81").Value = 26198.666
$ws.Range("M81").Value = -5638
$ws.Range("N81").Value = -28320.666
$ws.Range("H84").Value = 9199.4
$ws.Range("I84").Value = 3349.5
$ws.Range("J84").Value = 13099.333
$ws.Range("K84").Value = 33495
$ws.Range("L84").Value = 130993.33
$ws.Range("M84").Value = -28191
$ws.Range("N84").Value = -141601.33
$ws.Range("H126").Value = 156588.92
$ws.Range("I126").Value = 288700.16
$ws.Range("J126").Value = 2459.1667
$ws.Range("K126").Value = 866100.48
$ws.Range("L126").Value = 7377.500100000001
$ws.Range("M126").Value = -863630.48
$ws.Range("N126").Value = -12317.5001
$ws.Range("H132").Value = 30888.457
$ws.Range("I132").Value = 31698.941
$ws.Range("K132").Value = 95096.823
$ws.Range("M132").Value = -92566.823
